# Generate Report for Handoff
# Updates the localization-status report for the file
# "59e2ffa5-d6fd-4842-a146-a4f49be793e5.md" (row 3 in every sheet):
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The handoff timestamps are refreshed
#   - An error detail is recorded because the handback file is stale
#   - The "Error Detail" column is widened so the new message is readable

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 16:57:21"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-05 16:57:15"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/785e931c9ab35ca5d78293e724994887ecec00d6/e2e/59e2ffa5-d6fd-4842-a146-a4f49be793e5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5bde2b180a32b43bb76265a600d584ee41cb17b/e2e/59e2ffa5-d6fd-4842-a146-a4f49be793e5.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-05 16:57:21"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/785e931c9ab35ca5d78293e724994887ecec00d6/e2e/59e2ffa5-d6fd-4842-a146-a4f49be793e5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5bde2b180a32b43bb76265a600d584ee41cb17b/e2e/59e2ffa5-d6fd-4842-a146-a4f49be793e5.md."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664

Write-Output "Applied handoff report changes"
